$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$ws.Range("A9").Value = "test1"
$ws.Range("A10").Value = "test2"
$ws.Range("B9").Value = "Test 1"
$ws.Range("B10").Value = "Test 2"

$ws.Range("A11").Value = "testBodyCapsule"
$ws.Range("B11").Value = "Capsule"
$ws.Range("A12").Value = "testBodySphere"
$ws.Range("B12").Value = "Sphere"

$ws.Range("A4").Value = "none"
$ws.Range("B4").Value = "None"

$ws.Range("B6").Value = "Motility"
$ws.Range("B5").Value = "Body"

$ws.Range("A5").Value = "categoryBody"
$ws.Range("A6").Value = "categoryMotility"

$ws.Range("A7").Value = "nucleoid"
$ws.Range("B7").Value = "Nucleoid"

$ws.Range("A8").Value = "ribosome"
$ws.Range("B8").Value = "Ribosome"

$ws.Range("B8").Select()
